$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4251862888688249
$ws.Range("C2").Value = 0.9307327647082201

$ws.Range("B4").Value = 1.090838003737177
$ws.Range("C4").Value = 0.9461170227698321

$ws.Range("B6").Value = 1.191061667365153

$ws.Range("B7").Value = 1.271977756333025

$ws.Range("B8").Value = 1.73103996916734
$ws.Range("C8").Value = 0.9980171745939656

$ws.Range("B9").Value = 4.919241847624935
$ws.Range("C9").Value = 0.9411938629544762

$ws.Range("B10").Value = 0.8471066322618543
$ws.Range("C10").Value = 0.9974748013669915
